$wb = $excel.ActiveWorkbook

# --- Sheet "connections": re-lay rows 14-18. The old "CO2 Capture outflow
#     compressed CO2 -> CO2 Storage" row and the old "CO2 Storage inflow
#     electricity" row both slide down to make room for two new rows that
#     route the power process's waste heat into DRI's CO2 capture process.
$conn = $wb.Worksheets.Item("connections")

# Row 14 (was old row 15): CO2 Capture <- electricity
$conn.Range("B14").Value = "CO2 Capture"
$conn.Range("C14").Value = "simple_CO2capture"
$conn.Range("D14").Value = "inflow"
$conn.Range("E14").Value = "electricity"
$conn.Range("F14").Value = "electricity"
$conn.Range("G14").Value = "outflow"
$conn.Range("H14").Value = "simple_power"
$conn.Range("I14").Value = "power"
$conn.Range("F14").Font.Color = 0
$conn.Range("H14").NumberFormat = "@"

# Row 15 (new): power -> waste heat -> recovered heat -> CO2 Capture (replacing "heat")
$conn.Range("B15").Value = "power"
$conn.Range("C15").Value = "simple_power"
$conn.Range("D15").Value = "outflow"
$conn.Range("E15").Value = "waste heat"
$conn.Range("F15").Value = "recovered heat"
$conn.Range("G15").Value = "inflows"
$conn.Range("H15").Value = "simple_CO2capture"
$conn.Range("I15").Value = "CO2 Capture"
$conn.Range("J15").Value = "heat"
$conn.Range("C15").NumberFormat = "@"
$conn.Range("F15").Font.Color = 0
$conn.Range("H15").NumberFormat = "@"

# Row 16 (new): CO2 Capture <- heat (simple_heat outflow)
$conn.Range("B16").Value = "CO2 Capture"
$conn.Range("C16").Value = "simple_CO2capture"
$conn.Range("D16").Value = "inflow"
$conn.Range("E16").Value = "heat"
$conn.Range("F16").Value = "heat"
$conn.Range("G16").Value = "outflows"
$conn.Range("H16").Value = "simple_heat"
$conn.Range("I16").Value = "heat"
$conn.Range("C16").ClearFormats()
$conn.Range("F16").Font.Color = 0
$conn.Range("H16").NumberFormat = "@"

# Row 17 (was old row 14): CO2 Capture -> compressed CO2 -> CO2 Storage
$conn.Range("B17").Value = "CO2 Capture"
$conn.Range("C17").Value = "simple_CO2capture"
$conn.Range("D17").Value = "outflow"
$conn.Range("E17").Value = "compressed CO2"
$conn.Range("F17").Value = "compressed CO2"
$conn.Range("G17").Value = "inflows"
$conn.Range("H17").Value = "simple_CO2storage"
$conn.Range("I17").Value = "CO2 Storage"
$conn.Range("F17").Font.Color = 0
$conn.Range("H17").NumberFormat = "@"

# Row 18 (was old row 16): CO2 Storage <- electricity
$conn.Range("B18").Value = "CO2 Storage"
$conn.Range("C18").Value = "simple_CO2storage"
$conn.Range("D18").Value = "inflow"
$conn.Range("E18").Value = "electricity"
$conn.Range("F18").Value = "electricity"
$conn.Range("G18").Value = "outflow"
$conn.Range("H18").Value = "simple_power"
$conn.Range("I18").Value = "power"
$conn.Range("C18").NumberFormat = "@"
$conn.Range("F18").Font.Color = 0
$conn.Range("H18").NumberFormat = "@"

# --- Sheet "chains": register the new "heat" chain sourced from simple_heat.
$chains = $wb.Worksheets.Item("chains")
$chains.Range("A9").Value = "heat"
$chains.Range("B9").Value = "heat"
$chains.Range("C9").Value = "outflow"
$chains.Range("E9").Value = "simple_heat"
$chains.Range("E9").NumberFormat = "@"
